$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted before the existing row 800,
# pushing every following row down by one (800->801, ..., 902->903).
$ws.Rows.Item(800).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(800, 1).Value2 = 10
$ws.Cells.Item(800, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(800, 3).Value2 = "La Araucanía"
$ws.Cells.Item(800, 4).Value2 = 45077
$ws.Cells.Item(800, 5).Value2 = 9
$ws.Cells.Item(800, 6).Value2 = "Fruta"
$ws.Cells.Item(800, 7).Value2 = 100102
$ws.Cells.Item(800, 8).Value2 = "Cítricos"
$ws.Cells.Item(800, 9).Value2 = 100102004
$ws.Cells.Item(800, 10).Value2 = "Mandarina"
$ws.Cells.Item(800, 11).Value2 = "Clementina"
$ws.Cells.Item(800, 12).Value2 = "Primera"
$ws.Cells.Item(800, 13).Value2 = 150
$ws.Cells.Item(800, 14).Value2 = 15000
$ws.Cells.Item(800, 15).Value2 = 15000
$ws.Cells.Item(800, 16).Value2 = 15000
$ws.Cells.Item(800, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(800, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(800, 19).Value2 = 833
$ws.Cells.Item(800, 20).Value2 = 18
